$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 2 ("charging staff + freeze -> stop charging" / Bugfix row).
# This shifts all subsequent rows up by one, shrinks the table automatically,
# and updates dependent ranges (conditional formatting, sort state, etc).
$ws.Rows("2:2").Delete()
